$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the erroneous row (PROJECT_ID 152140018016610 / MEMBER_ID test) - row 7
$ws.Rows.Item(7).Delete()

# Update selection to match target (E6)
$ws.Range("E6").Select()
